$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated timestamp for the data refresh
$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 16:16"

# Full refreshed country table (sorted descending by "Casos totales"),
# País, Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes
$data = @(
    @("China", 81008, 41, 71740, 6013, 1927, 7, 3255),
    @("Italia", 47021, 0, 5129, 37860, 2655, 0, 4032),
    @("España", 24926, 3355, 2125, 21475, 939, 233, 1326),
    @("Alemania", 21652, 1804, 209, 21368, 2, 7, 75),
    @("Iran", 20610, 966, 7635, 11419, 0, 123, 1556),
    @("Estados Unidos", 19956, 573, 147, 19530, 64, 23, 279),
    @("Francia", 12612, 0, 1587, 10575, 1297, 0, 450),
    @("Corea del Sur", 8799, 147, 2612, 6085, 59, 8, 102),
    @("Suiza", 6186, 571, 15, 6113, 141, 2, 58),
    @("Reino Unido", 4072, 89, 65, 3830, 20, 0, 177),
    @("Paises Bajos", 3631, 637, 2, 3493, 210, 30, 136),
    @("Belgica", 2815, 558, 263, 2485, 288, 30, 67),
    @("Austria", 2814, 165, 9, 2797, 15, 2, 8),
    @("Noruega", 2040, 81, 1, 2032, 28, 0, 7),
    @("Suecia", 1764, 125, 16, 1728, 69, 4, 20),
    @("Dinamarca", 1326, 71, 1, 1312, 42, 4, 13),
    @("Portugal", 1280, 260, 5, 1263, 26, 6, 12),
    @("Malasia", 1183, 153, 114, 1061, 26, 5, 8),
    @("Canada", 1087, 0, 14, 1060, 1, 1, 13),
    @("Australia", 1072, 144, 46, 1019, 2, 0, 7),
    @("Japon", 1007, 0, 215, 757, 50, 0, 35),
    @("Brasil", 999, 29, 2, 985, 18, 1, 12),
    @("Chequia", 925, 92, 5, 920, 7, 0, 0),
    @("Israel", 883, 178, 36, 846, 15, 0, 1),
    @("Pakistan", 720, 219, 13, 704, 0, 0, 3),
    @("Crucero", 712, 0, 567, 137, 15, 0, 8),
    @("Irlanda", 683, 0, 5, 675, 6, 0, 3),
    @("Turquia", 670, 0, 6, 661, 3, 3, 9),
    @("Luxemburgo", 670, 186, 6, 656, 3, 0, 8),
    @("Chile", 537, 103, 8, 529, 7, 0, 0),
    @("Finlandia", 521, 71, 10, 510, 2, 1, 1),
    @("Grecia", 495, 0, 19, 463, 20, 3, 13),
    @("Islandia", 473, 64, 5, 468, 1, 0, 0),
    @("Catar", 470, 0, 10, 460, 6, 0, 0),
    @("Polonia", 452, 27, 13, 434, 3, 0, 5),
    @("Indonesia", 450, 81, 20, 392, 0, 6, 38),
    @("Singapur", 432, 47, 140, 290, 14, 2, 2),
    @("Ecuador", 426, 0, 3, 416, 2, 0, 7),
    @("Tailandia", 411, 89, 44, 366, 7, 0, 1),
    @("Arabia Saudita", 392, 48, 16, 376, 0, 0, 0),
    @("Eslovenia", 383, 42, 0, 382, 12, 0, 1),
    @("Rumania", 367, 59, 52, 315, 14, 0, 0),
    @("Filipinas", 307, 77, 13, 275, 1, 1, 19),
    @("Estonia", 306, 23, 2, 304, 1, 0, 0),
    @("Rusia", 306, 53, 16, 289, 0, 0, 1),
    @("Barein", 305, 7, 125, 179, 4, 0, 1),
    @("Egipto", 285, 0, 42, 235, 0, 0, 8),
    @("India", 283, 34, 23, 255, 0, 0, 5),
    @("Hong Kong", 273, 17, 98, 171, 4, 0, 4),
    @("Peru", 263, 0, 1, 258, 5, 0, 4),
    @("Sudafrica", 240, 38, 2, 238, 0, 0, 0),
    @("Irak", 214, 6, 51, 146, 0, 0, 17),
    @("Libano", 206, 29, 8, 194, 4, 0, 4),
    @("Mexico", 203, 39, 4, 197, 1, 1, 2),
    @("Panama", 200, 0, 1, 198, 7, 0, 1),
    @("Kuwait", 176, 17, 27, 149, 5, 0, 0),
    @("Croacia", 168, 38, 5, 162, 0, 0, 1),
    @("Bulgaria", 163, 36, 3, 157, 3, 0, 3),
    @("Armenia", 160, 24, 1, 159, 2, 0, 0),
    @("Colombia", 158, 13, 1, 157, 0, 0, 0),
    @("Argentina", 158, 0, 3, 151, 0, 1, 4),
    @("Taiwan", 153, 18, 28, 123, 0, 0, 2),
    @("Emiratos Arabes Unidos", 153, 13, 38, 113, 2, 0, 2),
    @("San Marino", 151, 0, 4, 127, 12, 6, 20),
    @("Serbia", 149, 14, 2, 146, 4, 0, 1),
    @("Eslovaquia", 138, 1, 7, 131, 2, 0, 0),
    @("Letonia", 124, 13, 1, 123, 0, 0, 0),
    @("Costa Rica", 113, 0, 2, 109, 2, 0, 2),
    @("Republica Dominicana", 112, 40, 0, 109, 0, 1, 3),
    @("Uruguay", 110, 0, 0, 110, 0, 0, 0),
    @("Hungria", 103, 18, 7, 92, 6, 0, 4),
    @("Argelia", 95, 1, 43, 40, 0, 1, 12),
    @("Vietnam", 94, 3, 17, 77, 2, 0, 0),
    @("Bosnia y Herzegovina", 93, 4, 2, 91, 1, 0, 0),
    @("Islas Feroe", 92, 12, 3, 89, 0, 0, 0),
    @("Principado de Andorra", 88, 13, 1, 87, 2, 0, 0),
    @("Marruecos", 86, 0, 2, 81, 1, 0, 3),
    @("Jordania", 85, 0, 1, 84, 0, 0, 0),
    @("Republica de Macedonia", 85, 9, 1, 84, 1, 0, 0),
    @("Brunei", 83, 5, 1, 82, 2, 0, 0),
    @("Lituania", 83, 14, 1, 81, 1, 0, 1),
    @("Sri Lanka", 77, 4, 3, 74, 2, 0, 0),
    @("Albania", 76, 6, 2, 72, 2, 0, 2),
    @("Bielorrusia", 76, 7, 15, 61, 0, 0, 0),
    @("Republica de Chipre", 75, 0, 0, 75, 1, 0, 0),
    @("Malta", 73, 9, 2, 71, 1, 0, 0),
    @("Moldavia", 66, 0, 1, 64, 3, 0, 1),
    @("Venezuela", 65, 0, 1, 64, 0, 0, 0),
    @("Burkina Faso", 64, 24, 5, 56, 0, 2, 3),
    @("Tunez", 60, 6, 1, 58, 7, 0, 1),
    @("Kazajistan", 53, 1, 0, 53, 0, 0, 0),
    @("Camboya", 53, 2, 1, 52, 0, 0, 0),
    @("Azerbaiyan", 53, 9, 11, 41, 0, 0, 1),
    @("Nueva Zelanda", 52, 0, 0, 52, 0, 0, 0),
    @("Oman", 52, 4, 13, 39, 0, 0, 0),
    @("Estado de Palestina", 52, 4, 17, 35, 0, 0, 0),
    @("Guadalupe", 51, 0, 0, 50, 4, 0, 1),
    @("Georgia", 49, 5, 1, 48, 1, 0, 0),
    @("Senegal", 47, 0, 5, 42, 0, 0, 0),
    @("Reunion", 45, 7, 0, 45, 0, 0, 0),
    @("Uzbekistan", 41, 8, 0, 41, 0, 0, 0),
    @("Ucrania", 41, 0, 1, 37, 0, 0, 3),
    @("Liechtenstein", 37, 9, 0, 37, 0, 0, 0),
    @("Martinica", 32, 0, 0, 31, 7, 0, 1),
    @("Camerun", 27, 0, 2, 25, 0, 0, 0),
    @("Honduras", 24, 0, 0, 24, 0, 0, 0),
    @("Afganistan", 24, 0, 1, 23, 0, 0, 0),
    @("Banglades", 24, 4, 3, 19, 0, 1, 2),
    @("Consejo Danes para los Refugiados", 23, 5, 0, 23, 0, 0, 0),
    @("Nigeria", 22, 10, 1, 21, 0, 0, 0),
    @("Cuba", 21, 0, 0, 20, 0, 0, 1),
    @("Ghana", 19, 3, 0, 19, 0, 0, 0),
    @("Bolivia", 19, 3, 0, 19, 0, 0, 0),
    @("Jamaica", 19, 0, 2, 16, 0, 0, 1),
    @("Paraguay", 18, 0, 0, 17, 1, 1, 1),
    @("Macao", 18, 1, 10, 8, 0, 0, 0),
    @("Puerto Rico", 17, 3, 0, 17, 0, 0, 0),
    @("Ruanda", 17, 0, 0, 17, 0, 0, 0),
    @("Togo", 16, 7, 0, 16, 0, 0, 0),
    @("Guayana Francesa", 15, 0, 0, 15, 0, 0, 0),
    @("Polinesia Francesa", 15, 4, 0, 15, 0, 0, 0),
    @("Guam", 15, 1, 0, 15, 0, 0, 0),
    @("Kirguistan", 14, 8, 0, 14, 0, 0, 0),
    @("Montenegro", 14, 0, 0, 14, 0, 0, 0),
    @("Mauricio", 14, 2, 0, 13, 0, 1, 1),
    @("Costa de Marfil", 14, 5, 1, 13, 0, 0, 0),
    @("Guatemala", 13, 1, 0, 12, 0, 0, 1),
    @("Maldivas", 13, 0, 3, 10, 0, 0, 0),
    @("Monaco", 11, 0, 0, 11, 0, 0, 0),
    @("Mongolia", 10, 4, 0, 10, 0, 0, 0),
    @("Gibraltar", 10, 0, 2, 8, 0, 0, 0),
    @("Trinidad yTobago", 9, 0, 0, 9, 0, 0, 0),
    @("Etiopia", 9, 0, 0, 9, 0, 0, 0),
    @("Seychelles", 7, 0, 0, 7, 0, 0, 0),
    @("Kenia", 7, 0, 0, 7, 0, 0, 0),
    @("Mayotte", 7, 0, 0, 7, 0, 0, 0),
    @("Guinea Ecuatorial", 6, 0, 0, 6, 0, 0, 0),
    @("Tanzania", 6, 0, 0, 6, 0, 0, 0),
    @("Barbados", 6, 0, 0, 6, 0, 0, 0),
    @("Islas Virgenes de los Estados Unidos", 6, 3, 0, 6, 0, 0, 0),
    @("Aruba", 5, 0, 1, 4, 0, 0, 0),
    @("Guyana", 5, 0, 0, 4, 0, 0, 1),
    @("Bahamas", 4, 0, 0, 4, 0, 0, 0),
    @("San Martin (Parte Francesa)", 4, 0, 0, 4, 0, 0, 0),
    @("Nueva Caledonia", 4, 2, 0, 4, 0, 0, 0),
    @("Surinam", 4, 0, 0, 4, 0, 0, 0),
    @("Gabon", 4, 0, 0, 3, 0, 0, 1),
    @("Zimbabue", 3, 2, 0, 3, 0, 0, 0),
    @("Liberia", 3, 1, 0, 3, 0, 0, 0),
    @("Congo", 3, 0, 0, 3, 0, 0, 0),
    @("El Salvador", 3, 2, 0, 3, 0, 0, 0),
    @("Republica de Africa Central", 3, 0, 0, 3, 0, 0, 0),
    @("Namibia", 3, 0, 0, 3, 0, 0, 0),
    @("Madagascar", 3, 0, 0, 3, 0, 0, 0),
    @("San Bartolome", 3, 0, 0, 3, 0, 0, 0),
    @("Cabo Verde", 3, 2, 0, 3, 0, 0, 0),
    @("Islas Caimanes", 3, 0, 0, 2, 0, 0, 1),
    @("Curazao", 3, 0, 0, 2, 0, 0, 1),
    @("Nicaragua", 2, 0, 0, 2, 0, 0, 0),
    @("Benin", 2, 0, 0, 2, 0, 0, 0),
    @("Guinea", 2, 0, 0, 2, 0, 0, 0),
    @("Mauritania", 2, 0, 0, 2, 0, 0, 0),
    @("Isla de Man", 2, 0, 0, 2, 0, 0, 0),
    @("Haiti", 2, 0, 0, 2, 0, 0, 0),
    @("Angola", 2, 1, 0, 2, 0, 0, 0),
    @("Butan", 2, 0, 0, 2, 0, 0, 0),
    @("Zambia", 2, 0, 0, 2, 0, 0, 0),
    @("Groenlandia", 2, 0, 0, 2, 0, 0, 0),
    @("Bermudas", 2, 0, 0, 2, 0, 0, 0),
    @("Fiyi", 2, 1, 0, 2, 0, 0, 0),
    @("Santa Lucia", 2, 0, 0, 2, 0, 0, 0),
    @("Sudan", 2, 0, 0, 1, 0, 0, 1),
    @("Suazilandia", 1, 0, 0, 1, 0, 0, 0),
    @("Montserrat", 1, 0, 0, 1, 0, 0, 0),
    @("Papua Nueva Guinea", 1, 0, 0, 1, 0, 0, 0),
    @("Santa Sede", 1, 0, 0, 1, 0, 0, 0),
    @("Antigua y Barbuda", 1, 0, 0, 1, 0, 0, 0),
    @("Timor Oriental", 1, 1, 0, 1, 0, 0, 0),
    @("Somalia", 1, 0, 0, 1, 0, 0, 0),
    @("Gambia", 1, 0, 0, 1, 0, 0, 0),
    @("San Vicente y las Granadinas", 1, 0, 0, 1, 0, 0, 0),
    @("San Martin (Parte Holandesa)", 1, 0, 0, 1, 0, 0, 0),
    @("Republica de Yibuti", 1, 0, 0, 1, 0, 0, 0),
    @("Niger", 1, 0, 0, 1, 0, 0, 0),
    @("Republica del Chad", 1, 0, 0, 1, 0, 0, 0),
    @("Nepal", 1, 0, 1, 0, 0, 0, 0)
)

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $row.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $row[$c]
    }
}

Write-Host "Updated $($data.Count) country rows"
